$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "판다스(Pandas) .groupby()로 할 수 있는 거의 모든 것! (통계량, 전처리)"
$ws.Range("E4").Value = "https://teddylee777.github.io/pandas/pandas-groupby"

$ws.Range("D36").Value = "Multimodal Representation Learning : How to narrow heterogeneity gap"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/373"

$ws.Range("D51").Value = "[jQuery] 특정 요소에 담긴 text 정보 얻기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/jQuery-%ED%8A%B9%EC%A0%95-%EC%9A%94%EC%86%8C%EC%97%90-%EB%8B%B4%EA%B8%B4-text-%EC%A0%95%EB%B3%B4-%EC%96%BB%EA%B8%B0"
